$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text so numeric-looking values
# like "0.9979" or "82.18" are not auto-converted to numbers by Excel.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D16","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.366.80'
$ws.Range("D3").Value = '1.844.91'
$ws.Range("D4").Value = '0.9979'
$ws.Range("D5").Value = '240.39'
$ws.Range("D6").Value = '0.6311'
$ws.Range("D7").Value = '0.9991'
$ws.Range("D8").Value = '0.07499'
$ws.Range("D9").Value = '0.2907'
$ws.Range("D10").Value = '24.41'
$ws.Range("D12").Value = '1.844.67'
$ws.Range("D13").Value = '5.007'
$ws.Range("D14").Value = '0.6805'
$ws.Range("D16").Value = '82.18'
$ws.Range("D17").Value = '2.104.47'
$ws.Range("D18").Value = '6.155'
$ws.Range("D19").Value = '29.380.37'
$ws.Range("D20").Value = '229.24'
$ws.Range("D22").Value = '0.9991'
$ws.Range("D23").Value = '7.450'
$ws.Range("D24").Value = '0.9985'
$ws.Range("D25").Value = '158.87'
$ws.Range("D26").Value = '0.1380'
$ws.Range("D27").Value = '8.413'
$ws.Range("D29").Value = '0.06389'
$ws.Range("D30").Value = '1.386'
$ws.Range("D31").Value = '1.472'
$ws.Range("D32").Value = '4.093'
$ws.Range("D33").Value = '4.053'
$ws.Range("D34").Value = '1.820'
$ws.Range("D35").Value = '1.142'
$ws.Range("D36").Value = '0.6987'
$ws.Range("D38").Value = '1.259.00'
$ws.Range("D39").Value = '2.833'
$ws.Range("D40").Value = '0.01821'
$ws.Range("D41").Value = '6.586'
$ws.Range("D42").Value = '0.9084'
$ws.Range("D43").Value = '0.9984'
$ws.Range("D44").Value = '2.005.88'
$ws.Range("D45").Value = '101.33'
$ws.Range("D46").Value = '66.32'
$ws.Range("D48").Value = '0.1178'
$ws.Range("D49").Value = '7.052'
$ws.Range("D50").Value = '1.705'
$ws.Range("D51").Value = '9.059'

# Other text columns (Coin name, Link, Volume %) are safe to set directly.
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("E17").Value = '  -3.77%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("E29").Value = '  +14.33%  '
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E35").Value = '  -1.80%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("E41").Value = '  +2.76%  '
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("E51").Value = '  +0.50%  '
